$d = $word.ActiveDocument

# 1. Update the "Curso (semestre ideal)" line: drop ", EM (8)"
$d.Content.Find.Execute(
    "Curso (semestre ideal): EF (7), EM (8)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Curso (semestre ideal): EF (7)", 2)

# 2. Remove the trailing "Requisitos" heading paragraph and its
#    "LOM3049 -  Termodinâmica de Máquinas  (Requisito)" bullet paragraph.
$count = $d.Paragraphs.Count
$reqHeading = $d.Paragraphs.Item($count - 1)
$reqBullet = $d.Paragraphs.Item($count)
$r = $d.Range($reqHeading.Range.Start, $reqBullet.Range.End)
$r.Delete()
